$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.237840333333334
$ws.Range("N2").Value = 12.713521
$ws.Range("O2").Value = 0.4475167411500002
$ws.Range("P2").Value = 0.4475167411500002
$ws.Range("Q2").Value = 114.5865650033952
$ws.Range("R2").Value = 1031.279085030557
$ws.Range("S2").Value = 0.03175663026368997
$ws.Range("T2").Value = 0.03175663026368998

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.231839666666667
$ws.Range("N3").Value = 15.695519
$ws.Range("O3").Value = 0.5524832588499998
$ws.Range("P3").Value = 0.5524832588499998
$ws.Range("Q3").Value = 141.4632192101248
$ws.Range("R3").Value = 1273.168972891123
$ws.Range("S3").Value = 0.0392052519266473
$ws.Range("T3").Value = 0.03920525192664731

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 345.566579
$ws.Range("H4").Value = 1036.699737
$ws.Range("I4").Value = 0.9069174311350353
$ws.Range("J4").Value = 0.9069174311350354
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.237840333333334
$ws.Range("N4").Value = 12.713521
$ws.Range("O4").Value = 0.4475167411500002
$ws.Range("P4").Value = 0.4475167411500002
$ws.Range("Q4").Value = 1464.45598633822
$ws.Range("R4").Value = 13180.10387704398
$ws.Range("S4").Value = 0.4058607332736807
$ws.Range("T4").Value = 0.4058607332736808

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.231839666666667
$ws.Range("N5").Value = 15.695519
$ws.Range("O5").Value = 0.5524832588499998
$ws.Range("P5").Value = 0.5524832588499998
$ws.Range("Q5").Value = 1807.9489354865
$ws.Range("R5").Value = 16271.5404193785
$ws.Range("S5").Value = 0.5010566978613546
$ws.Range("T5").Value = 0.5010566978613546

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.428738666666668
$ws.Range("H6").Value = 25.286216
$ws.Range("I6").Value = 0.0221206866746274
$ws.Range("J6").Value = 0.02212068667462741
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.237840333333334
$ws.Range("N6").Value = 12.713521
$ws.Range("O6").Value = 0.4475167411500002
$ws.Range("P6").Value = 0.4475167411500002
$ws.Range("Q6").Value = 35.71964868072623
$ws.Range("R6").Value = 321.4768381265361
$ws.Range("S6").Value = 0.009899377612629488
$ws.Range("T6").Value = 0.009899377612629491

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.428738666666668
$ws.Range("H7").Value = 25.286216
$ws.Range("I7").Value = 0.0221206866746274
$ws.Range("J7").Value = 0.02212068667462741
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.231839666666667
$ws.Range("N7").Value = 15.695519
$ws.Range("O7").Value = 0.5524832588499998
$ws.Range("P7").Value = 0.5524832588499998
$ws.Range("Q7").Value = 44.09780929623378
$ws.Range("R7").Value = 396.880283666104
$ws.Range("S7").Value = 0.01222130906199791
$ws.Range("T7").Value = 0.01222130906199792
